# Updated symbol list on Sun Dec 25 11:38:57 UTC 2022 with GitHub Actions
#
# Applies the price/coin-list refresh captured in the commit diff:
#   - refreshed Price (column D) values for many rows
#   - rows 41-43 rotated (KickToken/BKEXToken/CEJI reshuffled with new
#     coinranking links, rank labels in column E, and new prices)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new (text) value. Column D holds
# numeric-looking strings that Excel would otherwise silently coerce into
# real numbers, so every cell is forced to Text format before the write
# and restored to the default "Normal" style afterwards (matches the
# original workbook, where these are plain un-styled string cells).
$updates = @(
    @{ Cell = "D2";  Value = "244.81" }
    @{ Cell = "D3";  Value = "23.08" }
    @{ Cell = "D4";  Value = "5.434" }
    @{ Cell = "D6";  Value = "3.391" }
    @{ Cell = "D7";  Value = "0.8099" }
    @{ Cell = "D8";  Value = "0.9227" }
    @{ Cell = "D9";  Value = "0.1432" }
    @{ Cell = "D10"; Value = "0.07432" }
    @{ Cell = "D11"; Value = "0.03389" }
    @{ Cell = "D12"; Value = "0.03040" }
    @{ Cell = "D13"; Value = "0.09341" }
    @{ Cell = "D14"; Value = "3.955" }
    @{ Cell = "D15"; Value = "0.001594" }
    @{ Cell = "D16"; Value = "0.04828" }
    @{ Cell = "D18"; Value = "0.005424" }
    @{ Cell = "D19"; Value = "0.004157" }
    @{ Cell = "D20"; Value = "0.0009818" }
    @{ Cell = "D22"; Value = "3.661" }
    @{ Cell = "D23"; Value = "6.448" }
    @{ Cell = "D26"; Value = "0.1340" }
    @{ Cell = "D40"; Value = "0.03932" }

    @{ Cell = "B41"; Value = "KickToken" }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick" }
    @{ Cell = "D41"; Value = "0.006213" }
    @{ Cell = "E41"; Value = "40KickTokenKICK" }

    @{ Cell = "B42"; Value = "BKEXToken" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk" }
    @{ Cell = "D42"; Value = "0.1074" }
    @{ Cell = "E42"; Value = "41BKEXTokenBKK" }

    @{ Cell = "B43"; Value = "CEJI" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji" }
    @{ Cell = "D43"; Value = "0.002901" }
    @{ Cell = "E43"; Value = "42CEJICEJI" }

    @{ Cell = "D44"; Value = "0.007112" }
    @{ Cell = "D45"; Value = "0.00005135" }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = "Normal"
}
